$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column B (mirrors the user's interaction before editing)
$ws.Range("B:B").Select()

# Fix the "harvester" column: it was mistakenly filled with "Retrofitted_0684".
# Holly added "S.GISH" as the correct harvester value for bioSamples.
$ws.Range("B2:B19").Value = "S.GISH"
